$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "58.823.21"
$ws.Range("E2").Value = "  -3.39%  "

$ws.Range("D3").Value = "3.233.75"
$ws.Range("E3").Value = "  -3.73%  "

$ws.Range("E4").Value = "  +0.07%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "539.03"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -5.00%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "136.67"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -8.42%  "

$ws.Range("E7").Value = "  -0.02%  "

$ws.Range("D8").Value = "3.234.34"
$ws.Range("E8").Value = "  -3.71%  "

$ws.Range("E9").Value = "  -4.39%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.62"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.93%  "

$ws.Range("E11").Value = "  -5.78%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.395"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -4.32%  "

$ws.Range("D13").Value = "3.798.05"
$ws.Range("E13").Value = "  -3.53%  "

$ws.Range("E14").Value = "  -1.16%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "26.04"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -7.40%  "

$ws.Range("D16").Value = "3.243.30"
$ws.Range("E16").Value = "  -3.43%  "

$ws.Range("E17").Value = "  -6.14%  "

$ws.Range("D18").Value = "58.930.21"
$ws.Range("E18").Value = "  -3.31%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "5.92"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -7.15%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.28"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -6.09%  "

$ws.Range("E21").Value = "  -6.11%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "360.94"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.58%  "

$ws.Range("E23").Value = "  +0.12%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "70.41"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -6.67%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.520"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -7.17%  "

$ws.Range("D26").Value = "3.376.20"
$ws.Range("E26").Value = "  -3.71%  "

$ws.Range("E27").Value = "  -3.55%  "

$ws.Range("D28").Value = "0.0₃0969"
$ws.Range("E28").Value = "  -10.87%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.999"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.07%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.07"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -4.20%  "

$ws.Range("E31").Value = "  -0.02%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.93"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -6.98%  "

$ws.Range("E33").Value = "  -7.55%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "22.03"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.76%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.23"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -4.50%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "163.45"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.96%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.94"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -7.98%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.43"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -4.97%  "

$ws.Range("E39").Value = "  -6.50%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "26.40"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -9.65%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0710"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -5.55%  "

$ws.Range("D42").Value = "3.273.97"
$ws.Range("E42").Value = "  -3.50%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "41.12"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.78%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.717"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -5.49%  "

$ws.Range("B45").Value = "Filecoin"
$ws.Range("C45").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "4.03"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -6.26%  "

$ws.Range("B46").Value = "ONDO"
$ws.Range("C46").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.10"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.67%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.50"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -6.54%  "

$ws.Range("E48").Value = "  +0.07%  "

$ws.Range("D49").Value = "2.301.01"
$ws.Range("E49").Value = "  -7.96%  "

$ws.Range("E50").Value = "  -5.59%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "20.88"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -7.86%  "
